# Implemented getting number of lines for methods and classes.
# Adds two new worksheets (classNumberOfLines, methodNumberOfLines) with the
# "number of lines" metric for each class / method, matching the target
# OOXML diff. Also re-syncs the classFields sheet row order (a pure
# byproduct of the source tool re-running), so that the visible content of
# every sheet in the workbook matches the target state cell-for-cell.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Re-order the classFields rows (sheet "classFields") so the per-class
#    field groupings match the target state. The field set itself is
#    unchanged - only the row order within each class differs.
# ---------------------------------------------------------------------
$classFields = $wb.Worksheets.Item("classFields")

$fieldRows = @(
    @("pl.piomin.stock.StockComponentTests", "template", "private", "org.springframework.kafka.core.KafkaTemplate"),
    @("pl.piomin.stock.StockComponentTests", "product", "", "pl.piomin.stock.domain.Product"),
    @("pl.piomin.stock.StockComponentTests", "repository", "", "pl.piomin.stock.repository.ProductRepository"),
    @("pl.piomin.stock.StockComponentTests", "factory", "private", "org.springframework.kafka.core.ConsumerFactory"),
    @("pl.piomin.stock.StockComponentTests", "LOG", "private", "org.slf4j.Logger"),
    @("pl.piomin.stock.StockComponentTests", "kafka", "private", "org.springframework.kafka.test.EmbeddedKafkaBroker"),
    @("pl.piomin.stock.service.OrderManageService", "repository", "private", "pl.piomin.stock.repository.ProductRepository"),
    @("pl.piomin.stock.service.OrderManageService", "SOURCE", "private", "java.lang.String"),
    @("pl.piomin.stock.service.OrderManageService", "LOG", "private", "org.slf4j.Logger"),
    @("pl.piomin.stock.service.OrderManageService", "template", "private", "org.springframework.kafka.core.KafkaTemplate"),
    @("pl.piomin.stock.StockApp", "LOG", "private", "org.slf4j.Logger"),
    @("pl.piomin.stock.StockApp", "repository", "private", "pl.piomin.stock.repository.ProductRepository"),
    @("pl.piomin.stock.StockApp", "orderManageService", "", "pl.piomin.stock.service.OrderManageService")
)

$r = 6
foreach ($row in $fieldRows) {
    $classFields.Cells.Item($r, 1).Value = $row[0]
    $classFields.Cells.Item($r, 2).Value = $row[1]
    $classFields.Cells.Item($r, 3).Value = $row[2]
    $classFields.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Add the two new worksheets at the end of the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$classNumberOfLines = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$classNumberOfLines.Name = "classNumberOfLines"

$methodNumberOfLines = $wb.Worksheets.Add([System.Type]::Missing, $classNumberOfLines)
$methodNumberOfLines.Name = "methodNumberOfLines"

# ---------------------------------------------------------------------
# 3. Populate classNumberOfLines.
#    Column B holds numeric-looking text, so format it as Text first so
#    the values are written as strings (matching the source export,
#    which stores every value - numeric or not - as a shared string).
# ---------------------------------------------------------------------
$classData = @(
    @("Class Name", "Number of Lines"),
    @("pl.piomin.stock.domain.Product", "1"),
    @("pl.piomin.stock.StockAppTest", "6"),
    @("pl.piomin.stock.repository.ProductRepository", "3"),
    @("pl.piomin.stock.StockComponentTests", "43"),
    @("pl.piomin.stock.service.OrderManageService", "44"),
    @("pl.piomin.stock.KafkaContainerDevMode", "6"),
    @("pl.piomin.stock.StockApp", "29")
)

$classNumberOfLines.Range("B1:B8").NumberFormat = "@"

$r = 1
foreach ($row in $classData) {
    $classNumberOfLines.Cells.Item($r, 1).Value = $row[0]
    $classNumberOfLines.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 4. Populate methodNumberOfLines.
# ---------------------------------------------------------------------
$methodData = @(
    @("Class Name", "Method Signature", "Number of Lines"),
    @("pl.piomin.stock.domain.Product", "getId()", "1"),
    @("pl.piomin.stock.domain.Product", "setId(java.lang.Long)", "1"),
    @("pl.piomin.stock.domain.Product", "getName()", "1"),
    @("pl.piomin.stock.domain.Product", "setName(java.lang.String)", "1"),
    @("pl.piomin.stock.domain.Product", "getAvailableItems()", "1"),
    @("pl.piomin.stock.domain.Product", "setAvailableItems(int)", "1"),
    @("pl.piomin.stock.domain.Product", "getReservedItems()", "1"),
    @("pl.piomin.stock.domain.Product", "setReservedItems(int)", "1"),
    @("pl.piomin.stock.domain.Product", "toString()", "1"),
    @("pl.piomin.stock.StockAppTest", "main(java.lang.String[])", "3"),
    @("pl.piomin.stock.StockComponentTests", "eventAccept()", "1"),
    @("pl.piomin.stock.StockComponentTests", "eventReject()", "1"),
    @("pl.piomin.stock.StockComponentTests", "eventConfirm()", "1"),
    @("pl.piomin.stock.service.OrderManageService", "reserve(pl.piomin.base.domain.Order)", "3"),
    @("pl.piomin.stock.service.OrderManageService", "confirm(pl.piomin.base.domain.Order)", "3"),
    @("pl.piomin.stock.KafkaContainerDevMode", "kafka()", "3"),
    @("pl.piomin.stock.StockApp", "main(java.lang.String[])", "3"),
    @("pl.piomin.stock.StockApp", "onEvent(pl.piomin.base.domain.Order)", "2"),
    @("pl.piomin.stock.StockApp", "generateData()", "7")
)

$methodNumberOfLines.Range("C1:C20").NumberFormat = "@"

$r = 1
foreach ($row in $methodData) {
    $methodNumberOfLines.Cells.Item($r, 1).Value = $row[0]
    $methodNumberOfLines.Cells.Item($r, 2).Value = $row[1]
    $methodNumberOfLines.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$classNumberOfLines.Select()
